$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a per-month row index (0-based, reset at the start of each
# month). Convert it into one continuous running index for the whole sheet:
# for data rows 2..245, the new value is simply (row number - 2).
# Rows 2-17 already hold 0..15, so only rows 18-245 actually change value,
# but we simply (re)write the whole column for consistency/simplicity.

for ($r = 2; $r -le 245; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
